$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextCell "D2" '57.617.72'
Set-TextCell "E2" '  -4.39%  '

Set-TextCell "D3" '2.931.28'
Set-TextCell "E3" '  -2.44%  '

Set-TextCell "E4" '  +0.11%  '

Set-TextCell "D5" '549.20'
Set-TextCell "E5" '  -4.31%  '

Set-TextCell "D6" '130.38'
Set-TextCell "E6" '  +3.92%  '

Set-TextCell "E7" '  +0.10%  '

Set-TextCell "D8" '0.513'
Set-TextCell "E8" '  +1.95%  '

Set-TextCell "D9" '2.926.10'
Set-TextCell "E9" '  -2.47%  '

Set-TextCell "D10" '0.126'
Set-TextCell "E10" '  -4.21%  '

Set-TextCell "D11" '4.79'
Set-TextCell "E11" '  -5.50%  '

Set-TextCell "D12" '0.445'
Set-TextCell "E12" '  +0.99%  '

Set-TextCell "E13" '  -0.14%  '

Set-TextCell "D14" '32.90'
Set-TextCell "E14" '  +1.02%  '

Set-TextCell "E15" '  +0.21%  '

Set-TextCell "D16" '3.415.52'
Set-TextCell "E16" '  -2.26%  '

Set-TextCell "E17" '  +6.18%  '

Set-TextCell "D18" '2.931.08'
Set-TextCell "E18" '  -2.27%  '

Set-TextCell "D19" '57.629.73'
Set-TextCell "E19" '  -4.19%  '

Set-TextCell "D20" '416.78'
Set-TextCell "E20" '  -3.02%  '

Set-TextCell "D21" '13.16'
Set-TextCell "E21" '  +0.33%  '

Set-TextCell "E22" '  +2.40%  '

Set-TextCell "E23" '  -1.38%  '

Set-TextCell "D24" '13.02'
Set-TextCell "E24" '  +1.01%  '

Set-TextCell "D25" '79.81'
Set-TextCell "E25" '  +0.49%  '

Set-TextCell "E26" '  +0.00%  '

Set-TextCell "D27" '1.00'
Set-TextCell "E27" '  +0.13%  '

Set-TextCell "E28" '  -2.88%  '

Set-TextCell "D29" '7.47'
Set-TextCell "E29" '  +2.89%  '

Set-TextCell "E30" '  +1.14%  '

Set-TextCell "E31" '  -0.72%  '

Set-TextCell "D32" '5.98'
Set-TextCell "E32" '  -2.94%  '

Set-TextCell "D33" '0.0965'
Set-TextCell "E33" '  +2.30%  '

Set-TextCell "E34" '  +0.83%  '

Set-TextCell "D35" '0.936'
Set-TextCell "E35" '  -0.06%  '

Set-TextCell "E36" '  +0.57%  '

Set-TextCell "D37" '47.99'
Set-TextCell "E37" '  -4.72%  '

Set-TextCell "B38" 'Cosmos'
Set-TextCell "C38" 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell "D38" '8.70'
Set-TextCell "E38" '  +2.56%  '

Set-TextCell "B39" 'PEPE'
Set-TextCell "C39" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextCell "D39" '0.0₃0681'
Set-TextCell "E39" '  +1.35%  '

Set-TextCell "D40" '2.54'
Set-TextCell "E40" '  +2.92%  '

Set-TextCell "E41" '  -0.67%  '

Set-TextCell "D42" '377.21'
Set-TextCell "E42" '  +0.23%  '

Set-TextCell "E43" '  -2.89%  '

Set-TextCell "D44" '2.683.74'
Set-TextCell "E44" '  +0.33%  '

Set-TextCell "E45" '  +0.01%  '

Set-TextCell "E46" '  +1.36%  '

Set-TextCell "D47" '122.01'
Set-TextCell "E47" '  +1.05%  '

Set-TextCell "E48" '  +1.38%  '

Set-TextCell "E49" '  -1.88%  '

Set-TextCell "D50" '23.05'
Set-TextCell "E50" '  -2.36%  '

Set-TextCell "E51" '  -0.41%  '

